$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 715.4
$ws.Range("H10").Value = 9000
$ws.Range("J10").Value = 9000
$ws.Range("L10").Value = 9000
$ws.Range("N10").Value = -9586
$ws.Range("H15").Value = 226.59
$ws.Range("I15").Value = 226.59
$ws.Range("K15").Value = 679.77
$ws.Range("M15").Value = -510.77
$ws.Range("H18").Value = 1071.6875
$ws.Range("I18").Value = 724.7857
$ws.Range("J18").Value = 3500
$ws.Range("K18").Value = 724.7857
$ws.Range("L18").Value = 3500
$ws.Range("M18").Value = -440.7857
$ws.Range("N18").Value = -4068
$ws.Range("H129").Value = 4618.778
$ws.Range("J129").Value = 1052.6316
$ws.Range("L129").Value = 3157.8948
$ws.Range("N129").Value = -13157.8948
$ws.Range("H132").Value = 3681250
$ws.Range("I132").Value = 4242195
$ws.Range("K132").Value = 12726585
$ws.Range("M132").Value = -12724055
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("M4").ClearContents()
$ws.Range("H5").Value = 625.5
$ws.Range("I5").Value = 1000
$ws.Range("J5").Value = 500.66666
$ws.Range("K5").Value = 1000
$ws.Range("L5").Value = 500.66666
$ws.Range("M5").Value = -888
$ws.Range("N5").Value = -724.66666
$ws.Range("H9").Value = 22000
$ws.Range("J9").Value = 22000
$ws.Range("L9").Value = 22000
$ws.Range("N9").Value = -22340
$ws.Range("H20").Value = 22000
$ws.Range("J20").Value = 22000
$ws.Range("L20").Value = 22000
$ws.Range("N20").Value = -22540
$ws.Range("H23").Value = 53717.855
$ws.Range("I23").Value = 80006
$ws.Range("J23").Value = 34001.75
$ws.Range("K23").Value = 80006
$ws.Range("L23").Value = 34001.75
$ws.Range("M23").Value = -79747
$ws.Range("N23").Value = -34519.75
$ws.Range("H32").Value = 30111.666
$ws.Range("I32").Value = 8965.016
$ws.Range("J32").Value = 495338
$ws.Range("K32").Value = 8965.016
$ws.Range("L32").Value = 495338
$ws.Range("M32").Value = -8678.016
$ws.Range("N32").Value = -495912
$ws.Range("H37").Value = 0
$ws.Range("J37").Value = 0
$ws.Range("L37").ClearContents()
$ws.Range("N37").Value = 0
$ws.Range("H44").Value = 16000
$ws.Range("J44").Value = 16000
$ws.Range("L44").Value = 16000
$ws.Range("N44").Value = -16976
$ws.Range("H55").Value = 8418.888999999999
$ws.Range("J55").Value = 8346.25
$ws.Range("L55").Value = 8346.25
$ws.Range("N55").Value = -8976.25
$ws.Range("H63").Value = 3100
$ws.Range("I63").Value = 0
$ws.Range("J63").Value = 3100
$ws.Range("K63").Value = 0
$ws.Range("L63").ClearContents()
$ws.Range("M63").Value = 3100
$ws.Range("N63").Value = -4472
$ws.Range("H66").Value = 3100
$ws.Range("I66").Value = 0
$ws.Range("J66").Value = 3100
$ws.Range("K66").Value = 0
$ws.Range("L66").ClearContents()
$ws.Range("M66").Value = 15500
$ws.Range("N66").Value = -22364
$ws.Range("H110").Value = 100210430
$ws.Range("I110").Value = 125262850
$ws.Range("J110").Value = 750
$ws.Range("K110").Value = 125262850
$ws.Range("L110").Value = 750
$ws.Range("M110").Value = -125260805
$ws.Range("N110").Value = -4840
$ws.Range("H132").Value = 24898.73
$ws.Range("I132").Value = 32603.684
$ws.Range("K132").Value = 97811.052
$ws.Range("M132").Value = -95281.052
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 625.5
$ws.Range("I4").Value = 1000
$ws.Range("J4").Value = 500.66666
$ws.Range("K4").Value = 1000
$ws.Range("L4").Value = 500.66666
$ws.Range("M4").Value = -885
$ws.Range("N4").Value = -730.66666
$ws.Range("H15").Value = 3500
$ws.Range("J15").Value = 3500
$ws.Range("L15").Value = 3500
$ws.Range("N15").Value = -3954
$ws.Range("H22").Value = 270
$ws.Range("I22").Value = 248.33333
$ws.Range("J22").Value = 400
$ws.Range("K22").Value = 248.33333
$ws.Range("L22").Value = 400
$ws.Range("M22").Value = -75.33332999999999
$ws.Range("N22").Value = -746
$ws.Range("H59").Value = 68000
$ws.Range("J59").Value = 68000
$ws.Range("L59").Value = 68000
$ws.Range("N59").Value = -69694
$ws.Range("H82").Value = 18908.8
$ws.Range("I82").Value = 2063.75
$ws.Range("J82").Value = 30138.834
$ws.Range("K82").Value = 2063.75
$ws.Range("L82").Value = 30138.834
$ws.Range("M82").Value = -1680.75
$ws.Range("N82").Value = -30904.834
$ws.Range("H85").Value = 18908.8
$ws.Range("I85").Value = 2063.75
$ws.Range("J85").Value = 30138.834
$ws.Range("K85").Value = 2063.75
$ws.Range("L85").Value = 30138.834
$ws.Range("M85").Value = -737.75
$ws.Range("N85").Value = -32790.834
$ws.Range("H134").Value = 1796.3572
$ws.Range("I134").Value = 1676.6267
$ws.Range("J134").Value = 2794.111
$ws.Range("K134").Value = 5029.8801
$ws.Range("L134").Value = 8382.332999999999
$ws.Range("M134").Value = -2494.8801
$ws.Range("N134").Value = -13452.333
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H17").Value = 6196.6665
$ws.Range("I17").Value = 295
$ws.Range("K17").Value = 295
$ws.Range("M17").Value = -121
$ws.Range("H25").Value = 10250
$ws.Range("I25").Value = 500
$ws.Range("J25").Value = 20000
$ws.Range("K25").Value = 500
$ws.Range("L25").Value = 20000
$ws.Range("M25").Value = -326
$ws.Range("N25").Value = -20348
$ws.Range("H51").Value = 7944.75
$ws.Range("J51").Value = 7944.75
$ws.Range("L51").Value = 7944.75
$ws.Range("N51").Value = -9416.75
$ws.Range("H59").Value = 23612.8
$ws.Range("J59").Value = 25490
$ws.Range("L59").Value = 25490
$ws.Range("N59").Value = -27780
$ws.Range("H60").Value = 12996.667
$ws.Range("J60").Value = 14945.75
$ws.Range("L60").Value = 14945.75
$ws.Range("N60").Value = -15967.75
$ws.Range("H61").Value = 7944.75
$ws.Range("J61").Value = 7944.75
$ws.Range("L61").Value = 7944.75
$ws.Range("N61").Value = -8640.75
$ws.Range("H68").Value = 16360.158
$ws.Range("J68").Value = 16360.158
$ws.Range("L68").Value = 16360.158
$ws.Range("N68").Value = -17858.158
$ws.Range("H71").Value = 16360.158
$ws.Range("J71").Value = 16360.158
$ws.Range("L71").Value = 49080.474
$ws.Range("N71").Value = -56568.474
$ws.Range("H74").Value = 38058.43
$ws.Range("J74").Value = 38058.43
$ws.Range("L74").Value = 38058.43
$ws.Range("N74").Value = -39806.43
$ws.Range("H77").Value = 38058.43
$ws.Range("J77").Value = 38058.43
$ws.Range("L77").Value = 114175.29
$ws.Range("N77").Value = -122911.29
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 884.5625
$ws.Range("I113").Value = 1178.8
$ws.Range("J113").Value = 624.94116
$ws.Range("K113").Value = 3536.4
$ws.Range("L113").Value = 1874.82348
$ws.Range("M113").Value = -1366.4
$ws.Range("N113").Value = -6214.82348
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H18").Value = 14166.667
$ws.Range("J18").Value = 14166.667
$ws.Range("L18").Value = 14166.667
$ws.Range("N18").Value = -14752.667
$ws.Range("H43").Value = 1958.3667
$ws.Range("I43").Value = 863.0417
$ws.Range("J43").Value = 6339.6665
$ws.Range("K43").Value = 863.0417
$ws.Range("L43").Value = 6339.6665
$ws.Range("M43").Value = -712.0417
$ws.Range("N43").Value = -6641.6665
$ws.Range("H57").Value = 15800
$ws.Range("J57").Value = 17960
$ws.Range("L57").Value = 17960
$ws.Range("N57").Value = -19600
$ws.Range("H132").Value = 3663.8845
$ws.Range("I132").Value = 2666.6316
$ws.Range("J132").Value = 6370.7144
$ws.Range("K132").Value = 7999.8948
$ws.Range("L132").Value = 19112.1432
$ws.Range("M132").Value = -5469.8948
$ws.Range("N132").Value = -24172.1432
